# Update Cntf-Il6st.xlsx rows 2-13 with recomputed TPM-based statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ E = 3; F = 1; G = 0.7676553333333332; H = 2.302966; I = 0.3736977786965754; J = 0.3736977786965754; M = 73.202511; N = 219.607533; O = 0.3264904632507938; P = 0.3264904632507938; Q = 56.194297982542; R = 505.7486818428779; S = 0.1220087608824375; T = 0.1220087608824375 }
    3 = @{ E = 3; F = 1; G = 0.7676553333333332; H = 2.302966; I = 0.3736977786965754; J = 0.3736977786965754; M = 99.76728333333334; N = 299.30185; O = 0.4449719839907295; P = 0.4449719839907295; Q = 76.58688714301111; R = 689.2819842870999; S = 0.1662850419995437; T = 0.1662850419995437 }
    4 = @{ E = 3; F = 1; G = 0.7676553333333332; H = 2.302966; I = 0.3736977786965754; J = 0.3736977786965754; M = 39.54025133333334; N = 118.620754; O = 0.1763534446908907; P = 0.1763534446908907; Q = 30.35328481737378; R = 273.179563356364; S = 0.06590289054647523; T = 0.06590289054647522 }
    5 = @{ E = 3; F = 1; G = 0.7676553333333332; H = 2.302966; I = 0.3736977786965754; J = 0.3736977786965754; M = 11.70021233333333; N = 35.100637; O = 0.05218410806758597; P = 0.05218410806758598; Q = 8.981730398815776; R = 80.83557358934199; S = 0.01950108526811892; T = 0.01950108526811892 }
    6 = @{ E = 2; F = 0.6666666666666666; G = 0.811404; H = 2.434212; I = 0.3949948098567449; J = 0.3949948098567449; M = 73.202511; N = 219.607533; O = 0.3264904632507938; P = 0.3264904632507938; Q = 59.396810235444; R = 534.571292118996; S = 0.1289620384517879; T = 0.1289620384517879 }
    7 = @{ E = 2; F = 0.6666666666666666; G = 0.811404; H = 2.434212; I = 0.3949948098567449; J = 0.3949948098567449; M = 99.76728333333334; N = 299.30185; O = 0.4449719839907295; P = 0.4449719839907295; Q = 80.95157276580001; R = 728.5641548922; S = 0.1757616242079967; T = 0.1757616242079967 }
    8 = @{ E = 2; F = 0.6666666666666666; G = 0.811404; H = 2.434212; I = 0.3949948098567449; J = 0.3949948098567449; M = 39.54025133333334; N = 118.620754; O = 0.1763534446908907; P = 0.1763534446908907; Q = 32.08311809287201; R = 288.748062835848; S = 0.06965869535326036; T = 0.06965869535326036 }
    9 = @{ E = 2; F = 0.6666666666666666; G = 0.811404; H = 2.434212; I = 0.3949948098567449; J = 0.3949948098567449; M = 11.70021233333333; N = 35.100637; O = 0.05218410806758597; P = 0.05218410806758598; Q = 9.493599088116; R = 85.442391793044; S = 0.02061245184369995; T = 0.02061245184369995 }
    10 = @{ E = 3; F = 1; G = 0.475155; H = 1.425465; I = 0.2313074114466796; J = 0.2313074114466796; M = 73.202511; N = 219.607533; O = 0.3264904632507938; P = 0.3264904632507938; Q = 34.782539114205; R = 313.042852027845; S = 0.07551966391656838; T = 0.07551966391656839 }
    11 = @{ E = 3; F = 1; G = 0.475155; H = 1.425465; I = 0.2313074114466796; J = 0.2313074114466796; M = 99.76728333333334; N = 299.30185; O = 0.4449719839907295; P = 0.4449719839907295; Q = 47.40492351225; R = 426.64431161025; S = 0.102925317783189; T = 0.102925317783189 }
    12 = @{ E = 3; F = 1; G = 0.475155; H = 1.425465; I = 0.2313074114466796; J = 0.2313074114466796; M = 39.54025133333334; N = 118.620754; O = 0.1763534446908907; P = 0.1763534446908907; Q = 18.78774812229; R = 169.08973310061; S = 0.04079185879115511; T = 0.04079185879115511 }
    13 = @{ E = 3; F = 1; G = 0.475155; H = 1.425465; I = 0.2313074114466796; J = 0.2313074114466796; M = 11.70021233333333; N = 35.100637; O = 0.05218410806758597; P = 0.05218410806758598; Q = 5.559414391244999; R = 50.034729521205; S = 0.0120705709557671; T = 0.0120705709557671 }
}

$colIndex = @{ E = 5; F = 6; G = 7; H = 8; I = 9; J = 10; M = 13; N = 14; O = 15; P = 16; Q = 17; R = 18; S = 19; T = 20 }

foreach ($rowKey in $updates.Keys) {
    $rowNum = [int]$rowKey
    $rowData = $updates[$rowKey]
    foreach ($colLetter in $rowData.Keys) {
        $colNum = $colIndex[$colLetter]
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$colLetter]
    }
}

Write-Output "Updated $($updates.Count) rows"
